$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replicate the formatting of the last data row (row 4) onto the new row 5
$ws.Range("A4:AT4").Copy()
$ws.Range("A5:AT5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A5").Value = "2021年"

$values = @{
    "B5" = 75.119
    "C5" = 84.08199999999999
    "D5" = 55.484
    "E5" = 60.368
    "F5" = 60
    "G5" = 42.207
    "H5" = 57.712
    "I5" = 64.11799999999999
    "J5" = 63.58
    "K5" = 77.917
    "L5" = 52.544
    "M5" = 51.307
    "N5" = 41.915
    "O5" = 26.73
    "P5" = 52.452
    "Q5" = 54.165
    "R5" = 32.907
    "S5" = 31.942
    "T5" = 57.998
    "U5" = 26.337
    "V5" = 67.69799999999999
    "W5" = 84.127
    "X5" = 21.297
    "Y5" = 19.848
    "Z5" = 23.021
    "AA5" = 22.944
    "AB5" = 71.05200000000001
    "AC5" = 43.88
    "AD5" = 47.201
    "AE5" = 50.35
    "AF5" = 50.887
    "AG5" = 34.957
    "AH5" = 77.604
    "AI5" = 70.871
    "AJ5" = 47.908
    "AK5" = 51.341
    "AL5" = 23.99
    "AM5" = 46.238
    "AN5" = 54.916
    "AO5" = 66.334
    "AP5" = 41.209
    "AQ5" = 23.604
    "AR5" = 57.391
    "AS5" = 44.454
    "AT5" = 22.931
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
